# Auto-generated edit script: apply numeric corrections to leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 59000
$ws.Range("J109").Value = 59000
$ws.Range("L109").Value = 59000
$ws.Range("N109").Value = -61774
$ws.Range("H116").Value = 10604.17
$ws.Range("I116").Value = 9206.066000000001
$ws.Range("K116").Value = 9206.066000000001
$ws.Range("M116").Value = -5764.066000000001
$ws.Range("H132").Value = 32055.402
$ws.Range("I132").Value = 34081.42
$ws.Range("J132").Value = 8080.8335
$ws.Range("K132").Value = 102244.26
$ws.Range("L132").Value = 24242.5005
$ws.Range("M132").Value = -99714.25999999999
$ws.Range("N132").Value = -29302.5005
$ws.Range("H137").Value = 1850.1957
$ws.Range("I137").Value = 1129.0769
$ws.Range("J137").Value = 2787.65
$ws.Range("K137").Value = 3387.2307
$ws.Range("L137").Value = 8362.950000000001
$ws.Range("M137").Value = -837.2307000000001
$ws.Range("N137").Value = -13462.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1392.091
$ws.Range("I2").Value = 1238.0333
$ws.Range("K2").Value = 1238.0333
$ws.Range("M2").Value = -1125.0333
$ws.Range("H32").Value = 3755.34
$ws.Range("I32").Value = 1962.9176
$ws.Range("J32").Value = 13912.4
$ws.Range("K32").Value = 1962.9176
$ws.Range("L32").Value = 13912.4
$ws.Range("M32").Value = -1675.9176
$ws.Range("N32").Value = -14486.4
$ws.Range("H61").Value = 10244.286
$ws.Range("I61").Value = 12592.667
$ws.Range("J61").Value = 8483
$ws.Range("K61").Value = 12592.667
$ws.Range("L61").Value = 8483
$ws.Range("M61").Value = -12380.667
$ws.Range("N61").Value = -8907
$ws.Range("H116").Value = 1392.091
$ws.Range("I116").Value = 1238.0333
$ws.Range("K116").Value = 1238.0333
$ws.Range("M116").Value = 1055.9667
$ws.Range("H132").Value = 488531.94
$ws.Range("I132").Value = 667029.25
$ws.Range("K132").Value = 2001087.75
$ws.Range("M132").Value = -1998557.75
$ws.Range("H136").Value = 10244.286
$ws.Range("I136").Value = 12592.667
$ws.Range("J136").Value = 8483
$ws.Range("K136").Value = 37778.001
$ws.Range("L136").Value = 25449
$ws.Range("M136").Value = -35228.001
$ws.Range("N136").Value = -30549

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1392.091
$ws.Range("I3").Value = 1238.0333
$ws.Range("K3").Value = 1238.0333
$ws.Range("M3").Value = -1124.0333
$ws.Range("H35").Value = 68165.336
$ws.Range("J35").Value = 68165.336
$ws.Range("L35").Value = 68165.336
$ws.Range("N35").Value = -68785.336
$ws.Range("H94").Value = 346.7857
$ws.Range("I94").Value = 357.76315
$ws.Range("K94").Value = 357.76315
$ws.Range("M94").Value = 93.23685
$ws.Range("H105").Value = 4036.3547
$ws.Range("I105").Value = 3527.25
$ws.Range("J105").Value = 4579.4
$ws.Range("K105").Value = 3527.25
$ws.Range("L105").Value = 4579.4
$ws.Range("M105").Value = -1780.25
$ws.Range("N105").Value = -8073.4
$ws.Range("H107").Value = 2403.6128
$ws.Range("I107").Value = 2522.0715
$ws.Range("J107").Value = 1298
$ws.Range("K107").Value = 2522.0715
$ws.Range("L107").Value = 1298
$ws.Range("M107").Value = -602.0715
$ws.Range("N107").Value = -5138
$ws.Range("H134").Value = 734200.4399999999
$ws.Range("I134").Value = 859280
$ws.Range("J134").Value = 8739
$ws.Range("K134").Value = 2577840
$ws.Range("L134").Value = 26217
$ws.Range("M134").Value = -2575305
$ws.Range("N134").Value = -31287
$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9623.097
$ws.Range("I31").Value = 2774.4285
$ws.Range("K31").Value = 2774.4285
$ws.Range("M31").Value = -2479.4285
$ws.Range("H34").Value = 9623.097
$ws.Range("I34").Value = 2774.4285
$ws.Range("K34").Value = 2774.4285
$ws.Range("M34").Value = -2572.4285
$ws.Range("H105").Value = 1589.8889
$ws.Range("I105").Value = 1589.8889
$ws.Range("K105").Value = 1589.8889
$ws.Range("M105").Value = 157.1111000000001
$ws.Range("H134").Value = 2541.7
$ws.Range("I134").Value = 1724.8108
$ws.Range("K134").Value = 5174.4324
$ws.Range("M134").Value = -2639.4324
$ws.Range("H138").Value = 41734.5
$ws.Range("J138").Value = 43760
$ws.Range("L138").Value = 43760
$ws.Range("N138").Value = -54040

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 22.076923
$ws.Range("I12").Value = 12.5
$ws.Range("J12").Value = 26.333334
$ws.Range("K12").Value = 37.5
$ws.Range("L12").Value = 79.00000199999999
$ws.Range("M12").Value = 135.5
$ws.Range("N12").Value = -425.000002
$ws.Range("H33").Value = 399.33334
$ws.Range("I33").Value = 356.2
$ws.Range("J33").Value = 615
$ws.Range("K33").Value = 2137.2
$ws.Range("L33").Value = 3690
$ws.Range("M33").Value = -1854.2
$ws.Range("N33").Value = -4256
$ws.Range("H37").Value = 105999.4
$ws.Range("J37").Value = 105999.4
$ws.Range("L37").Value = 317998.2
$ws.Range("N37").Value = -318222.2
$ws.Range("H87").Value = 19371.908
$ws.Range("I87").Value = 16146.714
$ws.Range("K87").Value = 48440.142
$ws.Range("M87").Value = -47192.142
$ws.Range("H90").Value = 19371.908
$ws.Range("I90").Value = 16146.714
$ws.Range("K90").Value = 145320.426
$ws.Range("M90").Value = -139080.426
$ws.Range("H113").Value = 1825.0416
$ws.Range("I113").Value = 1220.8572
$ws.Range("J113").Value = 2073.8235
$ws.Range("K113").Value = 3662.5716
$ws.Range("L113").Value = 6221.470499999999
$ws.Range("M113").Value = -1492.5716
$ws.Range("N113").Value = -10561.4705
$ws.Range("H138").Value = 10383.5
$ws.Range("I138").Value = 10383.5
$ws.Range("K138").Value = 31150.5
$ws.Range("M138").Value = -26010.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 420.2
$ws.Range("I3").Value = 275.25
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 275.25
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -159.25
$ws.Range("N3").Value = -1232
$ws.Range("H10").Value = 1122.2222
$ws.Range("J10").Value = 1122.2222
$ws.Range("L10").Value = 1122.2222
$ws.Range("N10").Value = -1460.2222
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("H19").Value = 996
$ws.Range("J19").Value = 1200
$ws.Range("L19").Value = 1200
$ws.Range("H113").Value = 3026.4614
$ws.Range("J113").Value = 4197.2144
$ws.Range("L113").Value = 4197.2144
$ws.Range("N113").Value = -8537.214400000001
$ws.Range("H132").Value = 2063.4375
$ws.Range("I132").Value = 2042.7273
$ws.Range("J132").Value = 2109
$ws.Range("K132").Value = 6128.1819
$ws.Range("L132").Value = 6327
$ws.Range("M132").Value = -3598.1819
$ws.Range("N132").Value = -11387
$ws.Range("H134").Value = 43508.062
$ws.Range("J134").Value = 43508.062
$ws.Range("L134").Value = 130524.186
$ws.Range("N134").Value = -135594.186
$ws.Range("H139").Value = 152750
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 152750
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 152750
$ws.Range("N139").Value = -163030
$ws.Range("N19").Value = -1776
$ws.Range("N11").ClearContents()
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 1000
$ws.Range("I14").Value = 1000
$ws.Range("K14").Value = 1000
$ws.Range("H22").Value = 86411.414
$ws.Range("J22").Value = 3543.7
$ws.Range("L22").Value = 3543.7
$ws.Range("N22").Value = -4133.7
$ws.Range("H27").Value = 86411.414
$ws.Range("J27").Value = 3543.7
$ws.Range("L27").Value = 3543.7
$ws.Range("N27").Value = -3757.7
$ws.Range("M14").Value = -828

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 201113800
$ws.Range("I4").Value = 2759997.5
$ws.Range("J4").Value = 333349660
$ws.Range("K4").Value = 2759997.5
$ws.Range("L4").Value = 333349660
$ws.Range("M4").Value = -2759884.5
$ws.Range("N4").Value = -333349886
$ws.Range("H138").Value = 97997
$ws.Range("J138").Value = 97997
$ws.Range("L138").Value = 97997
$ws.Range("N138").Value = -108277
